$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 123..178 down by one (new row N gets old row N-1 data), working from
# the bottom up so source data is never clobbered before it is read.
for ($r = 179; $r -ge 124; $r--) {
    $src = $r - 1
    for ($c = 1; $c -le 18; $c++) {
        $ws.Cells.Item($r, $c).Value = $ws.Cells.Item($src, $c).Value2
    }
}

# Row 179 is a brand new row; give its date cell (column D) the same number
# format used throughout the table (the date number format).
$ws.Cells.Item(179, 4).NumberFormat = $ws.Cells.Item(178, 4).NumberFormat

# Row 123 gets brand new data
$ws.Range("D123").Value = 45119
$ws.Range("K123").Value = 7000
$ws.Range("L123").Value = 7000
$ws.Range("M123").Value = 7000
$ws.Range("N123").Value = "$/caja 50 unidades"
$ws.Range("O123").Value = "Región de Arica y Parinacota"
$ws.Range("P123").Value = 140
$ws.Range("Q123").Value = 50
